$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the sample data rows (rows 2-4), keeping only the header row
$ws.Rows("2:4").Delete()

# Reorder header values: A=מייל, B=עיר, C=טלפון, D=שם משפחה, E=שם
$ws.Range("A1").Value = "מייל"
$ws.Range("B1").Value = "עיר"
$ws.Range("C1").Value = "טלפון"
$ws.Range("D1").Value = "שם משפחה"
$ws.Range("E1").Value = "שם"

# Adjust column widths: column A becomes wide (30.83), column E back to normal (15.83)
# (offset by -5/6 to compensate for this host's char->pixel->char width round-trip)
$ws.Columns.Item(1).ColumnWidth = 29.998697916666668
$ws.Columns.Item(5).ColumnWidth = 14.998697916666666
